$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 131239963
$ws.Range("B4").Value = 57881
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "spel/sång"
$ws.Range("N4").Value = ""
$ws.Range("P4").Value = "Bjärnaryd, Sm"
$ws.Range("Q4").Value = 429548
$ws.Range("R4").Value = 6303017
$ws.Range("S4").Value = 25
$ws.Range("T4").Value = "Kronoberg"
$ws.Range("U4").Value = "Ljungby"
$ws.Range("V4").Value = "Småland"
$ws.Range("W4").Value = "Angelstad"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2026-02-20"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2026-02-20"
$ws.Range("AC4").Value = "Trummar flitigt."
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AT4").Value = ""
$ws.Range("AW4").Value = "Krister Wahlström"
$ws.Range("AX4").Value = "Krister Wahlström"
$ws.Range("AY4").Value = ""
